# Weekly update: two new price records for
# "Macroferia Regional de Talca" / Betarraga (Hortaliza) are inserted into
# the data block that starts at row 307, pushing the existing rows down by
# one position each time. The two rows that fall off the bottom of the
# original range (old rows 341 and 342) become the new rows 343 and 344,
# which is exactly what a plain row-insert produces.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceRow($Row, $Fecha, $Volumen, $PrecioMinimo, $PrecioMaximo, $PrecioPromedio, $PrecioKg, $Calidad) {
    $ws.Range("A$Row").Value = 5
    $ws.Range("B$Row").Value = "Macroferia Regional de Talca"
    $ws.Range("C$Row").Value = "Maule"
    $ws.Range("D$Row").Value = $Fecha
    $ws.Range("E$Row").Value = 7
    $ws.Range("F$Row").Value = 100114014
    $ws.Range("G$Row").Value = "Betarraga"
    $ws.Range("H$Row").Value = "Sin especificar"
    $ws.Range("I$Row").Value = $Calidad
    $ws.Range("J$Row").Value = $Volumen
    $ws.Range("K$Row").Value = $PrecioMinimo
    $ws.Range("L$Row").Value = $PrecioMaximo
    $ws.Range("M$Row").Value = $PrecioPromedio
    $ws.Range("N$Row").Value = "`$/paquete 5 unidades"
    $ws.Range("O$Row").Value = "Región del Maule"
    $ws.Range("P$Row").Value = $PrecioKg
    $ws.Range("Q$Row").Value = 5
    $ws.Range("R$Row").Value = "Hortaliza"
}

# Insert the first new weekly record above row 307; everything from the old
# row 307 down to the old row 342 shifts down by one row.
$ws.Rows.Item(307).Insert()
Set-PriceRow 307 44748 5000 700 700 700 140 "Primera"

# Insert the second new weekly record above what is now row 333 (the old
# row 332, after the first shift); everything below shifts down by one row
# again, so the old rows 341/342 end up as the new rows 343/344.
$ws.Rows.Item(333).Insert()
Set-PriceRow 333 44747 5000 700 700 700 140 "Primera"
